$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: change B3 from "dfds" to "a1234"
$ws.Range("B3").Value = "a1234"
# A3 keeps the same text "abc" but loses its special (Arial) formatting
$ws.Range("A3").ClearFormats()

# Row 4: change A4 from "visual_user" to "standard_user" (B4 "abcd" unchanged)
$ws.Range("A4").Value = "standard_user"

# New rows 6-9
$ws.Range("C6").Value = "Invalid"

$ws.Range("A7").Value = "standard_user"
$ws.Range("C7").Value = "Invalid"
# A7 picks up the same (Arial) formatting used by A2/A5
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B8").Value = "secret_sauce"
$ws.Range("C8").Value = "Invalid"

$ws.Range("A9").Value = 1234
$ws.Range("B9").Value = 1234
$ws.Range("C9").Value = "Invalid"

# Column A width change (closest representable value to the target 13.7109375)
$ws.Columns.Item(1).ColumnWidth = 12.83

# Update selection to C9
$ws.Range("C9").Select()
